$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the "Montevideo" / "New York" columns (Y and Z) ---
# Header row
$ws.Range("Y1").Value = "New York"
$ws.Range("Z1").Value = "Montevideo"

# Data rows 2-5 (swap Y/Z values)
$ws.Range("Y2").Value = 534
$ws.Range("Z2").Value = 552

$ws.Range("Y3").Value = 535
$ws.Range("Z3").Value = 553

$ws.Range("Y4").Value = 530
$ws.Range("Z4").Value = 552

$ws.Range("Y5").Value = 529
$ws.Range("Z5").Value = 551

# --- The old last row (5) had the "date-only" format; now that a new
# row (6) is appended below it, row 5 switches to full date-time format ---
$ws.Range("AH5").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- Append new row 6 with the next day's bunker prices ---
$ws.Range("A6").Value = 570
$ws.Range("B6").Value = 648
$ws.Range("C6").Value = 495
$ws.Range("D6").Value = 578
$ws.Range("E6").Value = 626
$ws.Range("F6").Value = 647
$ws.Range("G6").Value = 503
$ws.Range("H6").Value = 518
$ws.Range("I6").Value = 565
$ws.Range("J6").Value = 522
$ws.Range("K6").Value = 584
$ws.Range("L6").Value = 517
$ws.Range("M6").Value = 535
$ws.Range("N6").Value = 883
$ws.Range("O6").Value = 583
$ws.Range("P6").Value = 530
$ws.Range("Q6").Value = 503
$ws.Range("R6").Value = 528
$ws.Range("S6").Value = 608
$ws.Range("T6").Value = 646
$ws.Range("U6").Value = 577
$ws.Range("V6").Value = 490
$ws.Range("W6").Value = 550
$ws.Range("X6").Value = 515
$ws.Range("Y6").Value = 533
$ws.Range("Z6").Value = 553
$ws.Range("AA6").Value = 503
$ws.Range("AB6").Value = 545
$ws.Range("AC6").Value = 573.5
$ws.Range("AD6").Value = 520
$ws.Range("AE6").Value = 520
$ws.Range("AF6").Value = 527
$ws.Range("AG6").Value = 490
$ws.Range("AH6").Value = 45728
$ws.Range("AH6").NumberFormat = "YYYY-MM-DD"
$ws.Range("AI6").Value = 504
$ws.Range("AJ6").Value = 556
$ws.Range("AK6").Value = 522
$ws.Range("AL6").Value = 758
$ws.Range("AM6").Value = 649
$ws.Range("AN6").Value = 610
$ws.Range("AO6").Value = 503
$ws.Range("AP6").Value = 645
$ws.Range("AQ6").Value = 768
$ws.Range("AR6").Value = 515
$ws.Range("AS6").Value = 495
$ws.Range("AT6").Value = 567
$ws.Range("AU6").Value = 578
$ws.Range("AV6").Value = 639
